$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 276, shifting the existing row 276 (and all
# rows below it) down by one.
$ws.Rows(276).Insert()

# Populate the newly inserted row 276 with the new data point.
$ws.Cells.Item(276, 1).Value = 3
$ws.Cells.Item(276, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(276, 3).Value = "Coquimbo"
$ws.Cells.Item(276, 4).Value = 44706
$ws.Cells.Item(276, 5).Value = 5
$ws.Cells.Item(276, 6).Value = 100112001
$ws.Cells.Item(276, 7).Value = "Berenjena"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 133
$ws.Cells.Item(276, 11).Value = 6500
$ws.Cells.Item(276, 12).Value = 7000
$ws.Cells.Item(276, 13).Value = 6756
$ws.Cells.Item(276, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(276, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(276, 16).Value = 113
$ws.Cells.Item(276, 17).Value = 60
$ws.Cells.Item(276, 18).Value = "Hortaliza"
